$d = $word.ActiveDocument

# The "Resource(DoctorController" paragraph has its text split across two runs
# (an old stray "_GoBack" bookmark sits between them). Find that paragraph by
# its visible text so the edit is anchored to content, not a hard-coded index.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Resource(DoctorController*") {
        $target = $para
        break
    }
}

# Replace that paragraph with: the fixed-up "Resource(DoctorController)" line
# followed by the new "seeders / roles / middleware / schedules" notes block.
# The legacy "_GoBack" bookmark is recreated at the very end, right after the
# last new line, matching where Word left the user's cursor.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
  <w:p><w:r><w:t>Resource(DoctorController)</w:t></w:r></w:p>
  <w:p><w:r><w:t xml:space="preserve">3.- se hará uso de los seeders para crerar datos ficticios </w:t></w:r></w:p>
  <w:p><w:r><w:t>.- se refrescaran las tablas y se ejecutara el seeder</w:t></w:r></w:p>
  <w:p><w:r><w:t xml:space="preserve">Php artisan migrate:refresh </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:t>seed</w:t></w:r></w:p>
  <w:p></w:p>
  <w:p><w:r><w:t>/************************** tabla patient **************************/</w:t></w:r></w:p>
  <w:p></w:p>
  <w:p></w:p>
  <w:p><w:r><w:t>/******************** ROLES DE USUARIO **********************/</w:t></w:r></w:p>
  <w:p><w:r><w:t>MEDICOS = Gestionar horario, sus citas, pacientes(atendidos y por atender)</w:t></w:r></w:p>
  <w:p><w:r><w:lastRenderedPageBreak/><w:t>PACIENTE = (CITAS ACUALES, RESERVAR CITA NUEVA)</w:t></w:r></w:p>
  <w:p></w:p>
  <w:p><w:r><w:t>/**************** CREACION DE MIDDLEWARE ***********************/</w:t></w:r></w:p>
  <w:p><w:r><w:t>Php artisan make:middleware AdminMiddleware</w:t></w:r></w:p>
  <w:p></w:p>
  <w:p><w:r><w:t>/******************+ creación de schedules = horarios ***********************/</w:t></w:r></w:p>
  <w:p><w:r><w:t>Esta sección es especifica para doctores</w:t></w:r></w:p>
  <w:p><w:r><w:t>(día, activo, horainic, horafin, horaini, horafin, idMedico)</w:t></w:r></w:p>
  <w:p><w:r><w:t>// creo el modelo y la migración</w:t></w:r></w:p>
  <w:p><w:r><w:t>Php artisan make:model WorkDay -m</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
  <w:p></w:p>
</w:body>
</w:document>
'@

$target.Range.InsertXML($xml)
